$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains its original text formatting so that
# numeric-looking price strings (e.g. "63.086.34", "1.00") are not
# silently converted into numbers when we write new values.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "63.086.34"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "2.554.74"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "581.63"
$ws.Range("E5").Value = "  +1.48%  "
$ws.Range("D6").Value = "147.57"
$ws.Range("E6").Value = "  -2.26%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -1.19%  "
$ws.Range("D9").Value = "0.106"
$ws.Range("E9").Value = "  -0.88%  "
$ws.Range("D10").Value = "5.53"
$ws.Range("E10").Value = "  -3.95%  "
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("D12").Value = "0.355"
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("D13").Value = "27.28"
$ws.Range("E13").Value = "  -3.31%  "
$ws.Range("D14").Value = "3.011.22"
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("D15").Value = "63.001.41"
$ws.Range("E15").Value = "  -0.69%  "
$ws.Range("E16").Value = "  -0.79%  "
$ws.Range("D17").Value = "2.556.81"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").Value = "11.37"
$ws.Range("E18").Value = "  -2.72%  "
$ws.Range("D19").Value = "336.76"
$ws.Range("E19").Value = "  -1.73%  "
$ws.Range("E20").Value = "  -1.11%  "
$ws.Range("E21").Value = "  -2.43%  "
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").Value = "65.74"
$ws.Range("E23").Value = "  -0.69%  "
$ws.Range("B24").Value = "Kaspa"
$ws.Range("C24").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D24").Value = "0.170"
$ws.Range("E24").Value = "  -0.31%  "
$ws.Range("B25").Value = "Fetch.AI"
$ws.Range("C25").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D25").Value = "1.62"
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  -0.28%  "
$ws.Range("E28").Value = "  -1.35%  "
$ws.Range("E29").Value = "  +3.60%  "
$ws.Range("E30").Value = "  +2.75%  "
$ws.Range("D31").Value = "0.0₃0816"
$ws.Range("E31").Value = "  -2.72%  "
$ws.Range("D32").Value = "177.22"
$ws.Range("E32").Value = "  +0.26%  "
$ws.Range("D33").Value = "1.54"
$ws.Range("E33").Value = "  -1.35%  "
$ws.Range("D34").Value = "409.92"
$ws.Range("E34").Value = "  -2.59%  "
$ws.Range("D35").Value = "19.17"
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("D36").Value = "0.401"
$ws.Range("E36").Value = "  -1.82%  "
$ws.Range("D38").Value = "4.35"
$ws.Range("E38").Value = "  -2.39%  "
$ws.Range("D39").Value = "1.75"
$ws.Range("E39").Value = "  -0.46%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").Value = "39.74"
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("D42").Value = "151.34"
$ws.Range("E42").Value = "  -2.95%  "
$ws.Range("D43").Value = "3.77"
$ws.Range("E43").Value = "  -1.56%  "
$ws.Range("D44").Value = "21.01"
$ws.Range("E44").Value = "  -1.19%  "
$ws.Range("D45").Value = "0.0541"
$ws.Range("E45").Value = "  +1.13%  "
$ws.Range("E46").Value = "  -1.18%  "
$ws.Range("D47").Value = "0.0968"
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("D48").Value = "0.0240"
$ws.Range("E48").Value = "  +1.44%  "
$ws.Range("E49").Value = "  -2.21%  "
$ws.Range("D50").Value = "1.73"
$ws.Range("E50").Value = "  -5.65%  "
$ws.Range("D51").Value = "11.31"
$ws.Range("E51").Value = "  -0.07%  "
